$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorganize the footnote rows (34-38):
#  - B34: "Nota: ..." -> "Actualización: mayo 2024."
#  - B35: "FNI Fondo Nacional de Infraestructura." -> "Nota: Para Caminos propios, A partir de 2015 considera los accesos a los puentes nacionales."
#  - B36: "Fuente: CAPUFE. ..." -> "FNI Fondo Nacional de Infraestructura."
#  - B37: (new) -> "Fuente: CAPUFE. Caminos y Puentes Federales de Ingresos y Servicios Conexos."
#  - F37: "Ultima actualización: mayo 2024" -> cleared (style retained)
#  - F38: "Dirección General de Planeación" -> cleared (style retained)

$ws.Range("B34").Value = "Actualización: mayo 2024."
$ws.Range("B35").Value = "Nota: Para Caminos propios, A partir de 2015 considera los accesos a los puentes nacionales."
$ws.Range("B36").Value = "FNI Fondo Nacional de Infraestructura."
$ws.Range("B37").Value = "Fuente: CAPUFE. Caminos y Puentes Federales de Ingresos y Servicios Conexos."
$ws.Range("F37").ClearContents()
$ws.Range("F38").ClearContents()
